$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.2393064647912979, -0.8791878223419189, -0.1872301995754242),
    @(-0.1922698318958282, -0.9285151958465576, 0.8594874143600464),
    @(-3.570354700088501, -0.7802276611328125, -4.989242076873779),
    @(-1.221577763557434, 2.375196695327759, -2.503631114959717),
    @(2.165364503860474, 0.5566509366035461, -0.4453207552433014),
    @(0.2727513313293457, 0.5925393104553223, 0.4948008358478546),
    @(-0.2121229618787765, 1.80510675907135, 1.96942949295044),
    @(0.2755002379417419, 1.588096976280212, 2.037540912628174),
    @(0.8869763612747192, 0.8231409192085266, 1.362993121147156),
    @(0.2379320114850998, -0.7533495426177979, 0.1786780804395675),
    @(0.6151412725448608, 1.230893492698669, -0.3686571717262268),
    @(0.1403462886810302, 0.7915286421775818, 0.00137444678694),
    @(-0.3019201457500458, 0.041233405470848, -0.0345138870179653),
    @(-0.1716531366109848, -0.2657263875007629, -0.0448985956609249),
    @(-0.2078468948602676, 0.078801617026329, -0.026419922709465),
    @(0.0473420582711696, -0.0163406450301408, 0.0430659987032413),
    @(-0.0447458773851394, 0.1368338167667389, -0.0710130855441093),
    @(0.0578794814646244, 0.06536258012056349, -0.0146607663482427),
    @(-0.0494800843298435, 0.0968221426010131, 0.09865473955869671),
    @(-0.08109235763549801, 0.040775254368782, 0.0303905457258224)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
